# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row -> new F value for the "展览" sheet
$zhanlanUpdates = @{
    2  = 1600
    3  = 262
    4  = 8563
    6  = 68
    8  = 82
    9  = 1337
    10 = 110
    13 = 9212
    14 = 149
    17 = 170
    18 = 344
    19 = 6118
    20 = 1047
    21 = 63
    23 = 107
}

# Row -> new F value for the "全部类型" sheet
$quanbuUpdates = @{
    2  = 1600
    3  = 262
    4  = 8563
    6  = 68
    8  = 82
    9  = 1337
    10 = 110
    15 = 9212
    16 = 149
    19 = 170
    20 = 344
    21 = 6119
    22 = 1047
    23 = 63
    25 = 107
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Cells.Item($row, 6).Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Cells.Item($row, 6).Value = $quanbuUpdates[$row]
}
